$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row (row 1) strings:
#   *_old  -> *_FV2410   (columns A..J, i.e. 1..10)
#   *_new  -> *_FV2504   (columns L..U, i.e. 12..21)
# Column K ("diff") is left untouched.
$baseHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseHeaders[$i])_FV2410"
    $ws.Cells.Item(1, $i + 12).Value = "$($baseHeaders[$i])_FV2504"
}

# Turn the whole used range into an Excel table ("Table1") with the header row.
$tableRange = $ws.Range("A1:U81")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (split after row 1).
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
